## SBSquares_Task_Report.xlsx — add "Live Score Pulsing Squares" task rows
## plus the two blank section-separator rows that sit above "By Assignee"
## and "By Type" in the summary block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-create the two empty separator rows (108 and 121) ----------------
# These rows currently have no cells at all (they were skipped entirely in
# sheetData). Touching a row-level property with its own default value
# (OutlineLevel 0 is the default) is enough to make Excel emit an explicit
# empty <row r="..."/> element for it on save, without writing any cell
# data into it.
$ws.Rows.Item(108).OutlineLevel = 0
$ws.Rows.Item(121).OutlineLevel = 0

# --- Append the new task rows (159-165) -----------------------------------
$newTasks = @(
    @(159, "feature", "Live Score Pulsing Squares - DB migration (live_quarter_score JSONB column)", "architect", "Done"),
    @(160, "feature", "Live Score Pulsing Squares - CSS pulsing animations + ScoreBoard live display", "uiux-expert", "Done"),
    @(161, "feature", "Live Score Pulsing Squares - Grid/GridCell live winner/runner-up pulsing logic", "ui-dev-1", "Done"),
    @(162, "feature", "Live Score Pulsing Squares - Simulation system (API + fixture + admin runner)", "ui-dev-2", "Done"),
    @(163, "bugfix", "Live Score Pulsing Squares - Fix duplicate type in ScoreBoard, final review + build", "team-lead", "Done"),
    @(164, "bugfix", "Live pulsing: allow LIVE badge on squares that are already confirmed winners from earlier quarters", "Srini", "Done"),
    @(165, "bugfix", "Increase badge font sizes (winner/runner-up/LIVE/tentative) from 5-8px to 10px for mobile readability", "Srini", "Done")
)

foreach ($task in $newTasks) {
    $rowNum = $task[0]
    $ws.Cells.Item($rowNum, 1).Value = $rowNum
    $ws.Cells.Item($rowNum, 2).Value = $task[1]
    $ws.Cells.Item($rowNum, 3).Value = $task[2]
    $ws.Cells.Item($rowNum, 4).Value = $task[3]
    $ws.Cells.Item($rowNum, 5).Value = $task[4]
}
